$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.525.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5247'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3233'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7817'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07758'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.856.28'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.028'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007956'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.557.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.642'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.447'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.996'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '142.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.164'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.676'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08719'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.098'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.130'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.876'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7177'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.104'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.277'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.4849'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9010'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.939'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.689'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4168'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.998'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1232'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.8919'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.29%  '
